# Apply changes described by the diff:
# - E2 value changes from 0 to 50 (Eugenie Malayi's row)
# - Dependent formulas (F2, E9, F9) recalc automatically
# - Selection moves to I7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value for Eugenie Malayi (row 2), column E
$ws.Range("E2").Value = 50

# Force recalculation so formula cells (F2, E9, F9) pick up new cached values
$excel.Calculate()

# Move the active selection to I7, matching the updated sheet view state
$ws.Range("I7").Select()

$wb.Save()
